$wb = $excel.ActiveWorkbook

# --- Regular_Timetable ---
$ws = $wb.Worksheets.Item("Regular_Timetable")
$ws.Range("B2").Value  = "MINOR: Generative Ai [C102]"
$ws.Range("E6").Value  = "CS307 (Lab) [L107]"
$ws.Range("E7").Value  = "CS307 (Lab) [L107]"
$ws.Range("C8").Value  = "DA261 (Lab) [L207]"
$ws.Range("D8").Value  = "DA262 (Lab) [L106]"
$ws.Range("C9").Value  = "DA261 (Lab) [L207]"
$ws.Range("D9").Value  = "DA262 (Lab) [L106]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# --- PreMid_Timetable ---
$ws = $wb.Worksheets.Item("PreMid_Timetable")
$ws.Range("B2").Value  = "MINOR: Generative Ai [C102]"
$ws.Range("E6").Value  = "CS307 (Lab) [L107]"
$ws.Range("E7").Value  = "CS307 (Lab) [L107]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# --- PostMid_Timetable ---
$ws = $wb.Worksheets.Item("PostMid_Timetable")
$ws.Range("B2").Value  = "MINOR: Generative Ai [C102]"
$ws.Range("E6").Value  = "CS307 (Lab) [L107]"
$ws.Range("E7").Value  = "CS307 (Lab) [L107]"
$ws.Range("C8").Value  = "DA262 (Lab) [L207]"
$ws.Range("F8").Value  = "DA261 (Lab) [L207]"
$ws.Range("C9").Value  = "DA262 (Lab) [L207]"
$ws.Range("F9").Value  = "DA261 (Lab) [L207]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"
